# Apply the data-refresh edit described by the commit:
# "Update gh-pages to output generated at 456a3b4"
#
# This updates "want-to-go" counts (column F) for many rows across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets, flips one
# event's price cell to "已售罄" (sold out), and inserts one new duplicated
# row ("广州·Look Look动漫嘉年华") just above the "广州·第五届AP动漫嘉年华"
# row on sheet1 / sheet4 (pushing the rows below it down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 8398

$ws1.Range("F3").Value = 36637
$ws1.Range("G3").Value = "已售罄"

$ws1.Range("F5").Value = 617
$ws1.Range("F6").Value = 749
$ws1.Range("F7").Value = 471
$ws1.Range("F11").Value = 80
$ws1.Range("F13").Value = 495
$ws1.Range("F15").Value = 607
$ws1.Range("F17").Value = 450
$ws1.Range("F19").Value = 1139
$ws1.Range("F21").Value = 778
$ws1.Range("F22").Value = 2441
$ws1.Range("F23").Value = 933
$ws1.Range("F24").Value = 534
$ws1.Range("F26").Value = 1130

# Row 28 ("广州·Look Look动漫嘉年华") gets a refreshed "want to go" count.
$ws1.Range("F28").Value = 716

# Insert a new row before row 29 (the "广州·第五届AP动漫嘉年华" row),
# pushing it and everything below down by one. The new row duplicates the
# (now-updated) row 28 data, matching the upstream generator's output.
$ws1.Rows.Item(29).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown, [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromLeftOrAbove)

$ws1.Range("A28").Copy()
$ws1.Range("A29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("A29").Value = 28
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "2024-06-01"
$ws1.Range("C29").Value = "广州·Look Look动漫嘉年华"
$ws1.Range("D29").Value = "东沙大道16号 健康方舟"
$ws1.Range("E29").Value = "2024.06.01 10:00-06.02 17:30"
$ws1.Range("F29").Value = 716
$ws1.Range("G29").Value = 52.2
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$ws1.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png"

# Fix up the running index in column A for the rows that shifted down.
$ws1.Range("A30").Value = 29
$ws1.Range("A31").Value = 30

# Row 31 (was row 30, "广州·622排球少年only") also got a refreshed count.
$ws1.Range("F31").Value = 1125

# ---------------------------------------------------------------------
# Sheet "演出" (index 2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 53

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4) - mirrors the same edits as sheet 1/2, just at
# different row offsets because this sheet merges all three other sheets.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F3").Value = 8398

$ws4.Range("F5").Value = 36637
$ws4.Range("G5").Value = "已售罄"

$ws4.Range("F7").Value = 617
$ws4.Range("F8").Value = 749
$ws4.Range("F9").Value = 471
$ws4.Range("F17").Value = 80
$ws4.Range("F19").Value = 495
$ws4.Range("F20").Value = 53
$ws4.Range("F26").Value = 607
$ws4.Range("F28").Value = 450
$ws4.Range("F30").Value = 1139
$ws4.Range("F32").Value = 778
$ws4.Range("F33").Value = 2441
$ws4.Range("F34").Value = 933
$ws4.Range("F35").Value = 534
$ws4.Range("F37").Value = 1130

# Row 40 ("广州·Look Look动漫嘉年华") gets a refreshed "want to go" count.
$ws4.Range("F40").Value = 716

# Insert a new row before row 41 (the "广州·第五届AP动漫嘉年华" row),
# pushing it and everything below down by one, duplicating row 40's data.
$ws4.Rows.Item(41).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown, [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromLeftOrAbove)

$ws4.Range("A40").Copy()
$ws4.Range("A41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws4.Range("A41").Value = 40
$ws4.Range("B41").NumberFormat = "@"
$ws4.Range("B41").Value = "2024-06-01"
$ws4.Range("C41").Value = "广州·Look Look动漫嘉年华"
$ws4.Range("D41").Value = "东沙大道16号 健康方舟"
$ws4.Range("E41").Value = "2024.06.01 10:00-06.02 17:30"
$ws4.Range("F41").Value = 716
$ws4.Range("G41").Value = 52.2
$ws4.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$ws4.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png"

# Fix up the running index in column A for the rows that shifted down.
$ws4.Range("A42").Value = 41
$ws4.Range("A43").Value = 42

# Row 43 (was row 42, "广州·622排球少年only") also got a refreshed count.
$ws4.Range("F43").Value = 1125
